$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.653.33"
$ws.Range("E2").Value = "  +5.61%  "
$ws.Range("D3").Value = "2.257.53"
$ws.Range("E3").Value = "  +4.62%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "60.04"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "2.597.17"
$ws.Range("E13").Value = "  +4.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.18"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.41"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.835"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "2.264.45"
$ws.Range("E18").Value = "  +5.39%  "
$ws.Range("D19").Value = "41.622.66"
$ws.Range("E19").Value = "  +5.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "74.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").Value = "0.0₃0921"
$ws.Range("E21").Value = "  +8.30%  "
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +9.73%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.72"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "173.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("E31").Value = "  +2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.82"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.97%  "
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.93"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.24"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0637"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.88"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.45"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.88"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +12.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000234"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +49.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0238"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.78"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +14.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").Value = "1.515.55"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0944"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.54%  "
